$p = $ppt.ActivePresentation

# Add a new slide at the end using the "Title and Content" layout (same
# layout used by the other content slides in this deck, e.g. slide 6).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "For  restarting the mysql serever"

# --- Body / content placeholder --------------------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$line1 = " sudo /etc/init.d/mysql start"
$line2 = " sudo /etc/init.d/mysql stop"
$errMsg = "ERROR] Could not open file '/var/log/mysql/error.log' for error logging: Permission denied"
$line3 = "If it is failing with error : " + $errMsg + " "
$line4 = ""
$line5 = "https://support.plesk.com/hc/en-us/articles/115004039393-Unable-to-start-mysql-var-log-mysql-error-log-Permission-denied"
$line6 = ""

$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5 + "`r" + $line6

# Colour the error message (3rd paragraph) red.
$redStart = $line1.Length + 1 + $line2.Length + 1 + "If it is failing with error : ".Length + 1
$redLen = $errMsg.Length
$tr.Characters($redStart, $redLen).Font.Color.RGB = 255

# Turn the last-but-one paragraph's URL text into a hyperlink.
$linkStart = 1 + $line1.Length + 1 + $line2.Length + 1 + $line3.Length + 1 + $line4.Length + 1
$linkLen = $line5.Length
$linkRange = $tr.Characters($linkStart, $linkLen)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://support.plesk.com/hc/en-us/articles/115004039393-Unable-to-start-mysql-var-log-mysql-error-log-Permission-denied"

# Give the content placeholder a visible accent1-coloured outline, as in
# the other "how to fix it" slides.
$body.Line.Visible = $true
$body.Line.ForeColor.ObjectThemeColor = 5
